$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet 1: "Means"
# ----------------------------------------------------------------------------
$wsMeans = $wb.Worksheets.Item("Means")

# Header text updates (drop "Rural Areas (...)" prefix)
$wsMeans.Range("B1").Value = "National Average"
$wsMeans.Range("C1").Value = "State Average"

# Data updates (rows 2-10, columns B-G)
$meansData = @{
    2  = @(72, 83, 23, 35, 46, 33)
    3  = @(13, 9.4, 35, 29, 32, 57)
    4  = @(15, 7.3, 43, 36, 22, 11)
    5  = @(18, 6.9, 61, 49, 38, 20)
    6  = @(71, 62, 34, 39, 45, 47)
    7  = @(7.3, 7, 17, 14, 12, 11)
    8  = @(5.8, 6, 13, 13, 11, 10)
    9  = @(29, 23, 30, 30, 30, 30)
    10 = @(0.37, 0.3, 0.4, 0.37, 0.37, 0.37)
}

foreach ($row in $meansData.Keys) {
    $values = $meansData[$row]
    $wsMeans.Range("B$row").Value = $values[0]
    $wsMeans.Range("C$row").Value = $values[1]
    $wsMeans.Range("D$row").Value = $values[2]
    $wsMeans.Range("E$row").Value = $values[3]
    $wsMeans.Range("F$row").Value = $values[4]
    $wsMeans.Range("G$row").Value = $values[5]
}

# ----------------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ----------------------------------------------------------------------------
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Header text updates (drop "Rural Areas (...)" prefix)
$wsSD.Range("B1").Value = "National Average SD"
$wsSD.Range("C1").Value = "State Average SD"

# Data updates (rows 2-10, columns B-G)
$sdData = @{
    2  = @(27, 23, 13, 21, 27, 31)
    3  = @(23, 19, 22, 28, 30, 38)
    4  = @(16, 8.9, 17, 19, 18, 14)
    5  = @(22, 11, 22, 24, 23, 23)
    6  = @(37, 25, 11, 13, 15, 19)
    7  = @(8.7, 8.4, 9.4, 9.6, 9.7, 10)
    8  = @(7.8, 8.4, 10, 11, 10, 11)
    9  = @(10, 4.7, 0, 0, 0, 2.1)
    10 = @(0.14, 0.058, 0.022, 0.048, 0.048, 0.046)
}

foreach ($row in $sdData.Keys) {
    $values = $sdData[$row]
    $wsSD.Range("B$row").Value = $values[0]
    $wsSD.Range("C$row").Value = $values[1]
    $wsSD.Range("D$row").Value = $values[2]
    $wsSD.Range("E$row").Value = $values[3]
    $wsSD.Range("F$row").Value = $values[4]
    $wsSD.Range("G$row").Value = $values[5]
}
